$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old sample rows 4-11 entirely (only two demo rows remain afterwards).
$ws.Range("A4:A11").EntireRow.Delete()

# Row 2 -> new "fund 2" test investor (TSTF3 / Kotak / Domestic)
$ws.Range("A2").Value = "TSTF3"
$ws.Range("B2").Value = "TSTF3"
$ws.Range("C2").Value = [DateTime]"1991-10-18"
$ws.Range("D2").Value = "TSTFU2121D"
$ws.Range("H2").Value = "Domestic"
$ws.Range("I2").Value = "Kotak"

# Row 3 -> new "fund 2" test investor (TSTF4 / Axis / Foreign)
$ws.Range("A3").Value = "TSTF4"
$ws.Range("B3").Value = "TSTF4"
$ws.Range("C3").Value = [DateTime]"2000-05-30"
$ws.Range("D3").Value = "TSTFU2222D"
$ws.Range("H3").Value = "Foreign"
$ws.Range("I3").Value = "Axis"

# Drop the stray formatted-but-empty trailing cells that belonged to the
# removed rows' data (Q2/R2 and Q3), leaving only S2 behind.
$ws.Range("Q2").Clear()
$ws.Range("R2").Clear()
$ws.Range("Q3").Clear()

# Match the new active selection left behind by the edit.
$ws.Range("B4").Select()
